$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E20").Value = 8.4
$ws.Range("E21").Value = 5.8
$ws.Range("F21").Select()
